$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 70
$ws.Range("H70").Value = 1126.8422
$ws.Range("I70").Value = 1206.6875
$ws.Range("J70").Value = 701
$ws.Range("K70").Value = 3620.0625
$ws.Range("L70").Value = 2103
$ws.Range("M70").Value = -3350.0625
$ws.Range("N70").Value = -2643
# row 73
$ws.Range("H73").Value = 1126.8422
$ws.Range("I73").Value = 1206.6875
$ws.Range("J73").Value = 701
$ws.Range("K73").Value = 3620.0625
$ws.Range("L73").Value = 2103
$ws.Range("M73").Value = -2684.0625
$ws.Range("N73").Value = -3975
# row 103
$ws.Range("H103").Value = 600.36365
$ws.Range("I103").Value = 499.8
$ws.Range("J103").Value = 684.1667
$ws.Range("K103").Value = 1499.4
$ws.Range("L103").Value = 2052.5001
$ws.Range("M103").Value = -913.4000000000001
$ws.Range("N103").Value = -3224.5001
# row 135
$ws.Range("H135").Value = 3039.3333
$ws.Range("I135").Value = 271
$ws.Range("J135").Value = 8576
$ws.Range("K135").Value = 2439
$ws.Range("L135").Value = 77184
$ws.Range("M135").Value = 96
$ws.Range("N135").Value = -82254
# row 137
$ws.Range("H137").Value = 4847.227
$ws.Range("I137").Value = 4879.9443
$ws.Range("J137").Value = 4700
$ws.Range("K137").Value = 14639.8329
$ws.Range("L137").Value = 14100
$ws.Range("M137").Value = -12089.8329
$ws.Range("N137").Value = -19200
# row 138
$ws.Range("H138").Value = 5155.568
$ws.Range("I138").Value = 3326.3462
$ws.Range("J138").Value = 6020.291
$ws.Range("K138").Value = 9979.0386
$ws.Range("L138").Value = 18060.873
$ws.Range("M138").Value = -4839.0386
$ws.Range("N138").Value = -28340.873

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 32546.31
$ws.Range("I32").Value = 16554.355
$ws.Range("J32").Value = 83187.5
$ws.Range("K32").Value = 16554.355
$ws.Range("L32").Value = 83187.5
$ws.Range("M32").Value = -16267.355
$ws.Range("N32").Value = -83761.5
# row 43
$ws.Range("H43").Value = 12044.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12044.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12044.5
$ws.Range("N43").Value = -12670.5
# row 61
$ws.Range("H61").Value = 3151.8667
$ws.Range("I61").Value = 2136.6086
$ws.Range("J61").Value = 6487.7144
$ws.Range("K61").Value = 2136.6086
$ws.Range("L61").Value = 6487.7144
$ws.Range("M61").Value = -1924.6086
$ws.Range("N61").Value = -6911.7144
# row 74
$ws.Range("H74").Value = 1416.5454
$ws.Range("I74").Value = 931.5789
$ws.Range("J74").Value = 4488
$ws.Range("K74").Value = 931.5789
$ws.Range("L74").Value = 4488
$ws.Range("M74").Value = -57.57889999999998
$ws.Range("N74").Value = -6236
# row 77
$ws.Range("H77").Value = 1416.5454
$ws.Range("I77").Value = 931.5789
$ws.Range("J77").Value = 4488
$ws.Range("K77").Value = 4657.8945
$ws.Range("L77").Value = 22440
$ws.Range("M77").Value = -289.8945000000003
$ws.Range("N77").Value = -31176
# row 132
$ws.Range("H132").Value = 2955.7083
$ws.Range("I132").Value = 2611.6667
$ws.Range("J132").Value = 3987.8333
$ws.Range("K132").Value = 7835.000100000001
$ws.Range("L132").Value = 11963.4999
$ws.Range("M132").Value = -5305.000100000001
$ws.Range("N132").Value = -17023.4999
# row 136
$ws.Range("H136").Value = 3151.8667
$ws.Range("I136").Value = 2136.6086
$ws.Range("J136").Value = 6487.7144
$ws.Range("K136").Value = 6409.825800000001
$ws.Range("L136").Value = 19463.1432
$ws.Range("M136").Value = -3859.825800000001
$ws.Range("N136").Value = -24563.1432

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 107
$ws.Range("H107").Value = 2260.4
$ws.Range("I107").Value = 1952.1666
$ws.Range("J107").Value = 3493.3333
$ws.Range("K107").Value = 1952.1666
$ws.Range("L107").Value = 3493.3333
$ws.Range("M107").Value = -32.16660000000002
$ws.Range("N107").Value = -7333.3333
# row 134
$ws.Range("H134").Value = 1053.9429
$ws.Range("I134").Value = 1041.4117
$ws.Range("J134").Value = 1480
$ws.Range("K134").Value = 3124.2351
$ws.Range("L134").Value = 4440
$ws.Range("M134").Value = -589.2351000000003
$ws.Range("N134").Value = -9510

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 4253.2856
$ws.Range("I31").Value = 4171.125
$ws.Range("J31").Value = 4362.8335
$ws.Range("K31").Value = 4171.125
$ws.Range("L31").Value = 4362.8335
$ws.Range("M31").Value = -3876.125
$ws.Range("N31").Value = -4952.8335
# row 34
$ws.Range("H34").Value = 4253.2856
$ws.Range("I34").Value = 4171.125
$ws.Range("J34").Value = 4362.8335
$ws.Range("K34").Value = 4171.125
$ws.Range("L34").Value = 4362.8335
$ws.Range("M34").Value = -3969.125
$ws.Range("N34").Value = -4766.8335
# row 58
$ws.Range("H58").Value = 1762.5264
$ws.Range("I58").Value = 1831.5172
$ws.Range("J58").Value = 1540.2222
$ws.Range("K58").Value = 1831.5172
$ws.Range("L58").Value = 1540.2222
$ws.Range("M58").Value = -1628.5172
$ws.Range("N58").Value = -1946.2222
# row 62
$ws.Range("H62").Value = 2300
$ws.Range("I62").Value = 2100
$ws.Range("J62").Value = 2600
$ws.Range("K62").Value = 2100
$ws.Range("L62").Value = 2600
$ws.Range("M62").Value = -1476
$ws.Range("N62").Value = -3848
# row 65
$ws.Range("H65").Value = 2300
$ws.Range("I65").Value = 2100
$ws.Range("J65").Value = 2600
$ws.Range("K65").Value = 10500
$ws.Range("L65").Value = 13000
$ws.Range("M65").Value = -7380
$ws.Range("N65").Value = -19240
# row 132
$ws.Range("H132").Value = 1433.1063
$ws.Range("I132").Value = 1274.0625
$ws.Range("J132").Value = 1772.4
$ws.Range("K132").Value = 3822.1875
$ws.Range("L132").Value = 5317.200000000001
$ws.Range("M132").Value = -1292.1875
$ws.Range("N132").Value = -10377.2
# row 136
$ws.Range("H136").Value = 1762.5264
$ws.Range("I136").Value = 1831.5172
$ws.Range("J136").Value = 1540.2222
$ws.Range("K136").Value = 5494.5516
$ws.Range("L136").Value = 4620.6666
$ws.Range("M136").Value = -2944.5516
$ws.Range("N136").Value = -9720.6666

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 34
$ws.Range("H34").Value = 1322.5518
$ws.Range("I34").Value = 1036
$ws.Range("J34").Value = 1382.25
$ws.Range("K34").Value = 3108
$ws.Range("L34").Value = 4146.75
$ws.Range("M34").Value = -3024
$ws.Range("N34").Value = -4314.75
# row 39
$ws.Range("H39").Value = 1107.2
$ws.Range("I39").Value = 887.5
$ws.Range("J39").Value = 1149.0476
$ws.Range("K39").Value = 2662.5
$ws.Range("L39").Value = 3447.142800000001
$ws.Range("M39").Value = -2368.5
$ws.Range("N39").Value = -4035.142800000001
# row 113
$ws.Range("H113").Value = 1569.8125
$ws.Range("I113").Value = 7945
$ws.Range("J113").Value = 659.0714
$ws.Range("K113").Value = 23835
$ws.Range("L113").Value = 1977.2142
$ws.Range("M113").Value = -21665
$ws.Range("N113").Value = -6317.2142
# row 141
$ws.Range("H141").Value = 2886.652
$ws.Range("I141").Value = 2699.125
$ws.Range("J141").Value = 2986.6667
$ws.Range("K141").Value = 8097.375
$ws.Range("L141").Value = 8960.000100000001
$ws.Range("M141").Value = -2917.375
$ws.Range("N141").Value = -19320.0001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 2325
$ws.Range("I132").Value = 1911.8823
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 5735.6469
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -3205.6469
$ws.Range("N132").Value = -19058
# row 134
$ws.Range("H134").Value = 21236
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 21236
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 63708
$ws.Range("N134").Value = -68778
# row 135
$ws.Range("H135").Value = 77028.57000000001
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 77028.57000000001
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 77028.57000000001
$ws.Range("N135").Value = -87168.57000000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 55
$ws.Range("H55").Value = 399.03705
$ws.Range("I55").Value = 320.52173
$ws.Range("J55").Value = 850.5
$ws.Range("K55").Value = 320.52173
$ws.Range("L55").Value = 850.5
$ws.Range("M55").Value = -147.52173
$ws.Range("N55").Value = -1196.5
# row 68
$ws.Range("H68").Value = 1761.5393
$ws.Range("I68").Value = 875.62
$ws.Range("J68").Value = 2897.3333
$ws.Range("K68").Value = 875.62
$ws.Range("L68").Value = 2897.3333
$ws.Range("M68").Value = -126.62
$ws.Range("N68").Value = -4395.3333
# row 71
$ws.Range("H71").Value = 1761.5393
$ws.Range("I71").Value = 875.62
$ws.Range("J71").Value = 2897.3333
$ws.Range("K71").Value = 4378.1
$ws.Range("L71").Value = 14486.6665
$ws.Range("M71").Value = -634.1000000000004
$ws.Range("N71").Value = -21974.6665
# row 136
$ws.Range("H136").Value = 3202.0889
$ws.Range("I136").Value = 2687.25
$ws.Range("J136").Value = 3790.476
$ws.Range("K136").Value = 8061.75
$ws.Range("L136").Value = 11371.428
$ws.Range("M136").Value = -5511.75
$ws.Range("N136").Value = -16471.428

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 3102
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3102
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3102
$ws.Range("N62").Value = -4350
# row 64
$ws.Range("H64").Value = 14980
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 14980
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 14980
$ws.Range("N64").Value = -15476
# row 65
$ws.Range("H65").Value = 3102
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3102
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15510
$ws.Range("N65").Value = -21750
# row 67
$ws.Range("H67").Value = 14980
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 14980
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 14980
$ws.Range("N67").Value = -16696
# row 136
$ws.Range("H136").Value = 957.75
$ws.Range("I136").Value = 899.7826
$ws.Range("J136").Value = 1224.4
$ws.Range("K136").Value = 2699.3478
$ws.Range("L136").Value = 3673.2
$ws.Range("M136").Value = -149.3478
$ws.Range("N136").Value = -8773.200000000001
# row 137
$ws.Range("H137").Value = 41970
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 41970
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 41970
$ws.Range("N137").Value = -52170
